$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsWVR = $wb.Worksheets.Item("WVR")

# --- ALC ---
# Row 5
$wsALC.Range("H5").Value = 67.2
$wsALC.Range("I5").Value = 34
$wsALC.Range("J5").Value = 200
$wsALC.Range("K5").Value = 34
$wsALC.Range("L5").Value = 200
$wsALC.Range("M5").Value = 81
$wsALC.Range("N5").Value = -430

# Row 100
$wsALC.Range("H100").Value = 6667942.5
$wsALC.Range("I100").Value = 10417782
$wsALC.Range("J100").Value = 1560.6666
$wsALC.Range("K100").Value = 10417782
$wsALC.Range("L100").Value = 1560.6666
$wsALC.Range("M100").Value = -10417241
$wsALC.Range("N100").Value = -2642.6666

# Row 103
$wsALC.Range("H103").Value = 383050.34
$wsALC.Range("I103").Value = 686.8889
$wsALC.Range("J103").Value = 813209.25
$wsALC.Range("K103").Value = 2060.6667
$wsALC.Range("L103").Value = 2439627.75
$wsALC.Range("M103").Value = -1474.6667
$wsALC.Range("N103").Value = -2440799.75

# Row 118
$wsALC.Range("H118").Value = 673.86664
$wsALC.Range("I118").Value = 542.3333
$wsALC.Range("J118").Value = 1200
$wsALC.Range("K118").Value = 1626.9999
$wsALC.Range("L118").Value = 3600
$wsALC.Range("M118").Value = 30.00009999999997
$wsALC.Range("N118").Value = -6914

# Row 123
$wsALC.Range("H123").Value = 32550.4
$wsALC.Range("J123").Value = 32550.4
$wsALC.Range("L123").Value = 32550.4
$wsALC.Range("N123").Value = -42350.4

# Row 132
$wsALC.Range("H132").Value = 1517.5
$wsALC.Range("I132").Value = 1153.2759
$wsALC.Range("J132").Value = 5038.3335
$wsALC.Range("K132").Value = 3459.8277
$wsALC.Range("L132").Value = 15115.0005
$wsALC.Range("M132").Value = -929.8277000000003
$wsALC.Range("N132").Value = -20175.0005

# --- ARM ---
# Row 61
$wsARM.Range("H61").Value = 1610.5294
$wsARM.Range("I61").Value = 1552.3077
$wsARM.Range("J61").Value = 1799.75
$wsARM.Range("K61").Value = 1552.3077
$wsARM.Range("L61").Value = 1799.75
$wsARM.Range("M61").Value = -1340.3077
$wsARM.Range("N61").Value = -2223.75

# Row 86
$wsARM.Range("H86").Value = 40157
$wsARM.Range("J86").Value = 40157
$wsARM.Range("L86").Value = 40157
$wsARM.Range("N86").Value = -42529

# Row 89
$wsARM.Range("H89").Value = 40157
$wsARM.Range("J89").Value = 40157
$wsARM.Range("L89").Value = 120471
$wsARM.Range("N89").Value = -132327

# Row 123
$wsARM.Range("H123").Value = 30420
$wsARM.Range("J123").Value = 30420
$wsARM.Range("L123").Value = 30420
$wsARM.Range("N123").Value = -40220

# Row 127
$wsARM.Range("H127").Value = 65353.332
$wsARM.Range("J127").Value = 65353.332
$wsARM.Range("L127").Value = 65353.332
$wsARM.Range("N127").Value = -75273.33199999999

# Row 132
$wsARM.Range("H132").Value = 4171.963
$wsARM.Range("I132").Value = 1233.2046
$wsARM.Range("J132").Value = 17102.5
$wsARM.Range("K132").Value = 3699.6138
$wsARM.Range("L132").Value = 51307.5
$wsARM.Range("M132").Value = -1169.6138
$wsARM.Range("N132").Value = -56367.5

# Row 136
$wsARM.Range("H136").Value = 1610.5294
$wsARM.Range("I136").Value = 1552.3077
$wsARM.Range("J136").Value = 1799.75
$wsARM.Range("K136").Value = 4656.9231
$wsARM.Range("L136").Value = 5399.25
$wsARM.Range("M136").Value = -2106.9231
$wsARM.Range("N136").Value = -10499.25

# --- BSM ---
# Row 63
$wsBSM.Range("H63").Value = 35600
$wsBSM.Range("J63").Value = 35600
$wsBSM.Range("L63").Value = 35600
$wsBSM.Range("N63").Value = -36972

# Row 66
$wsBSM.Range("H66").Value = 35600
$wsBSM.Range("J66").Value = 35600
$wsBSM.Range("L66").Value = 106800
$wsBSM.Range("N66").Value = -113664

# Row 68
$wsBSM.Range("H68").Value = 42095
$wsBSM.Range("J68").Value = 42095
$wsBSM.Range("L68").Value = 42095
$wsBSM.Range("N68").Value = -43717

# Row 71
$wsBSM.Range("H71").Value = 42095
$wsBSM.Range("J71").Value = 42095
$wsBSM.Range("L71").Value = 126285
$wsBSM.Range("N71").Value = -134397

# Row 75
$wsBSM.Range("H75").Value = 11738
$wsBSM.Range("I75").Value = 11738
$wsBSM.Range("J75").Value = 0
$wsBSM.Range("K75").Value = 11738
$wsBSM.Range("L75").Value = 0
$wsBSM.Range("M75").Value = -10802
$wsBSM.Range("N75").ClearContents()

# Row 78
$wsBSM.Range("H78").Value = 11738
$wsBSM.Range("I78").Value = 11738
$wsBSM.Range("J78").Value = 0
$wsBSM.Range("K78").Value = 35214
$wsBSM.Range("L78").Value = 0
$wsBSM.Range("M78").Value = -30534
$wsBSM.Range("N78").ClearContents()

# Row 82
$wsBSM.Range("H82").Value = 26076.77
$wsBSM.Range("I82").Value = 13333
$wsBSM.Range("J82").Value = 29899.9
$wsBSM.Range("K82").Value = 13333
$wsBSM.Range("L82").Value = 29899.9
$wsBSM.Range("M82").Value = -12950
$wsBSM.Range("N82").Value = -30665.9

# Row 85
$wsBSM.Range("H85").Value = 26076.77
$wsBSM.Range("I85").Value = 13333
$wsBSM.Range("J85").Value = 29899.9
$wsBSM.Range("K85").Value = 13333
$wsBSM.Range("L85").Value = 29899.9
$wsBSM.Range("M85").Value = -12007
$wsBSM.Range("N85").Value = -32551.9

# Row 87
$wsBSM.Range("H87").Value = 50000
$wsBSM.Range("J87").Value = 50000
$wsBSM.Range("L87").Value = 50000
$wsBSM.Range("N87").Value = -52496

# Row 90
$wsBSM.Range("H90").Value = 50000
$wsBSM.Range("J90").Value = 50000
$wsBSM.Range("L90").Value = 150000
$wsBSM.Range("N90").Value = -162480

# --- CUL ---
# Row 132
$wsCUL.Range("H132").Value = 1738066.9
$wsCUL.Range("I132").Value = 2151.1428
$wsCUL.Range("J132").Value = 1951249.5
$wsCUL.Range("K132").Value = 19360.2852
$wsCUL.Range("L132").Value = 17561245.5
$wsCUL.Range("M132").Value = -16830.2852
$wsCUL.Range("N132").Value = -17566305.5

# --- GSM ---
# Row 80
$wsGSM.Range("H80").Value = 2758
$wsGSM.Range("I80").Value = 2770
$wsGSM.Range("J80").Value = 2747.0908
$wsGSM.Range("K80").Value = 2770
$wsGSM.Range("L80").Value = 2747.0908
$wsGSM.Range("M80").Value = -1772
$wsGSM.Range("N80").Value = -4743.0908

# Row 83
$wsGSM.Range("H83").Value = 2758
$wsGSM.Range("I83").Value = 2770
$wsGSM.Range("J83").Value = 2747.0908
$wsGSM.Range("K83").Value = 13850
$wsGSM.Range("L83").Value = 13735.454
$wsGSM.Range("M83").Value = -8858
$wsGSM.Range("N83").Value = -23719.454

# Row 107
$wsGSM.Range("H107").Value = 1440.2174
$wsGSM.Range("I107").Value = 827.2222
$wsGSM.Range("J107").Value = 1834.2858
$wsGSM.Range("K107").Value = 827.2222
$wsGSM.Range("L107").Value = 1834.2858
$wsGSM.Range("M107").Value = 1092.7778
$wsGSM.Range("N107").Value = -5674.2858

# Row 132
$wsGSM.Range("H132").Value = 2736.4866
$wsGSM.Range("I132").Value = 1862.2941
$wsGSM.Range("J132").Value = 3479.55
$wsGSM.Range("K132").Value = 5586.8823
$wsGSM.Range("L132").Value = 10438.65
$wsGSM.Range("M132").Value = -3056.8823
$wsGSM.Range("N132").Value = -15498.65

# --- WVR ---
# Row 56
$wsWVR.Range("H56").Value = 45976
$wsWVR.Range("I56").Value = 0
$wsWVR.Range("J56").Value = 45976
$wsWVR.Range("K56").Value = 0
$wsWVR.Range("L56").Value = 45976
$wsWVR.Range("M56").ClearContents()
$wsWVR.Range("N56").Value = -47404

# Row 64
$wsWVR.Range("H64").Value = 25111.2
$wsWVR.Range("J64").Value = 25111.2
$wsWVR.Range("L64").Value = 25111.2
$wsWVR.Range("N64").Value = -25607.2

# Row 67
$wsWVR.Range("H67").Value = 25111.2
$wsWVR.Range("J67").Value = 25111.2
$wsWVR.Range("L67").Value = 25111.2
$wsWVR.Range("N67").Value = -26827.2

# Row 123
$wsWVR.Range("H123").Value = 29885.4
$wsWVR.Range("J123").Value = 29885.4
$wsWVR.Range("L123").Value = 29885.4
$wsWVR.Range("N123").Value = -39685.4

# Row 124
$wsWVR.Range("H124").Value = 61333.332
$wsWVR.Range("J124").Value = 61333.332
$wsWVR.Range("L124").Value = 61333.332
$wsWVR.Range("N124").Value = -71153.33199999999
